# Post-Test Survey Responses.xlsx - apply author's edits
# 1) Re-word / re-number the 12 question headers in row 1 (B,D,F,H,J,L,N,P,R,S,U,W)
# 2) Shrink the custom-width column range from A:AO (41 cols) down to A:AC (29 cols)
# 3) Move the active selection from Y15 to AK23

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header text updates -------------------------------------------------
$ws.Range("B1").Value = "Q1. Did you find it easy to find the required buttons? (1 for very hard and 10 for very easy)"
$ws.Range("D1").Value = "Q2.  Did you recognize different icons and labels? (1 for very hard and 10 for very easy)"
$ws.Range("F1").Value = "Q3. Text and colors were easy for the eyes? (1 for uneasy and 10 for very clear)"
$ws.Range("H1").Value = "Q4. The application is easy to use with minimum instructions from the tester. (1 for very hard and 10 for very easy)"
$ws.Range("J1").Value = "Q5. Steps are easy to follow, and buttons are easy to remember  (1 for very hard and 10 for very easy)"
$ws.Range("L1").Value = "Q6. You managed to finish the tasks as fast as possible. (1 for slow and 10 for very fast)"
$ws.Range("N1").Value = "Q7. Everything was clear and easy to preform without frustration (1 for not clear and hard and 10 for very clear and easy)"
$ws.Range("P1").Value = "Q8. You managed to perform well in labs without difficulties in using the VR equipment  (1 for difficult and 10 for very easy)"
$ws.Range("R1").Value = "Q9. You managed to perform well in labs by hearing ROBO instructions and following them.  (1 for instructions were not clear and 10 for  instructions were very clear)"
$ws.Range("S1").Value = "Q10. ROBO  instructions were clear and easy to follow  (1 for hard and 10 for very easy)   "
$ws.Range("U1").Value = "Q11. You found the application interesting. (1 for not interesting and 10 for very interesting)"
$ws.Range("W1").Value = "Q12. You found the gamification functions fun and encouraged competitive behavior. . (1 for not interested and 10 for very encouraging)"

# --- 2) Drop the custom 18.85546875 width from columns AD:AO (30-41) -------
# so the sheet's <cols> collapses back to a single min=1 max=29 run.
$ws.Columns("AD:AO").Delete()
$ws.Columns("AD:AO").Insert()

# --- 3) Update the saved selection -----------------------------------------
$ws.Range("AK23").Select()
